$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 24 ("「経験」الخبرة ..."), shifting subsequent rows up.
$ws.Rows("24:24").Delete()
